$wb = $excel.ActiveWorkbook

# --- Sheet "Résumé" (tournament summary) ---
$ws1 = $wb.Worksheets.Item("Résumé")

# Type : DOUBLE -> SIMPLE
$ws1.Range("B2").Value = "SIMPLE"

# Nombre de place : 20 -> 15
$ws1.Range("B3").Value = 15

# Nombre de participants : 13 -> 0
$ws1.Range("B4").Value = 0

# Fin des inscriptions : 2023-05-30 -> 2023-07-08 01:59:00.0
$ws1.Range("B5").Value = "2023-07-08 01:59:00.0"

# Début du tournois : 2023-06-01 -> 2023-07-09 11:00:00.0
$ws1.Range("B6").Value = "2023-07-09 11:00:00.0"

# Fin du tournois : 2023-06-02 -> 2023-07-09 20:00:00.0
$ws1.Range("B7").Value = "2023-07-09 20:00:00.0"

# Widen column B to fit the longer date/time strings
$ws1.Columns("B").ColumnWidth = 19.88671875

# --- Sheet "Equipe" (team list) ---
$ws2 = $wb.Worksheets.Item("Equipe")

# Remove all registered team rows, keep only the header row
$ws2.Rows("2:13").Delete()
